$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.030.97'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '3.564.57'
$ws.Range("E3").Value = '  -2.30%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '192.58'
$ws.Range("E5").Value = '  +1.72%  '

# Row 6
$ws.Range("D6").Value = '571.27'
$ws.Range("E6").Value = '  -3.51%  '

# Row 7
$ws.Range("D7").Value = '3.562.91'
$ws.Range("E7").Value = '  -2.06%  '

# Row 8
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  -0.19%  '

# Row 9
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").Value = '0.673'
$ws.Range("E10").Value = '  -3.22%  '

# Row 11
$ws.Range("D11").Value = '0.149'
$ws.Range("E11").Value = '  -1.88%  '

# Row 12
$ws.Range("D12").Value = '55.07'
$ws.Range("E12").Value = '  -3.72%  '

# Row 13
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +0.35%  '

# Row 14
$ws.Range("D14").Value = '9.82'
$ws.Range("E14").Value = '  -3.02%  '

# Row 15
$ws.Range("D15").Value = '4.144.41'
$ws.Range("E15").Value = '  -2.10%  '

# Row 16
$ws.Range("D16").Value = '3.564.55'
$ws.Range("E16").Value = '  -2.33%  '

# Row 17
$ws.Range("E17").Value = '  -1.15%  '

# Row 18
$ws.Range("D18").Value = '66.967.70'
$ws.Range("E18").Value = '  -0.73%  '

# Row 19
$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  -1.45%  '

# Row 20
$ws.Range("D20").Value = '18.19'
$ws.Range("E20").Value = '  -2.86%  '

# Row 21
$ws.Range("E21").Value = '  -4.04%  '

# Row 22
$ws.Range("D22").Value = '401.79'
$ws.Range("E22").Value = '  +1.24%  '

# Row 23
$ws.Range("D23").Value = '4.15'
$ws.Range("E23").Value = '  -4.72%  '

# Row 24
$ws.Range("D24").Value = '12.21'
$ws.Range("E24").Value = '  +9.87%  '

# Row 25
$ws.Range("D25").Value = '85.41'
$ws.Range("E25").Value = '  -1.50%  '

# Row 26
$ws.Range("D26").Value = '2.90'
$ws.Range("E26").Value = '  -1.02%  '

# Row 27
$ws.Range("D27").Value = '12.52'
$ws.Range("E27").Value = '  +1.32%  '

# Row 28
$ws.Range("E28").Value = '  +1.09%  '

# Row 29
$ws.Range("D29").Value = '3.72'
$ws.Range("E29").Value = '  +2.42%  '

# Row 30
$ws.Range("D30").Value = '7.82'
$ws.Range("E30").Value = '  +6.85%  '

# Row 31
$ws.Range("D31").Value = '8.95'
$ws.Range("E31").Value = '  -2.67%  '

# Row 32
$ws.Range("D32").Value = '31.08'
$ws.Range("E32").Value = '  -1.85%  '

# Row 33
$ws.Range("D33").Value = '648.89'
$ws.Range("E33").Value = '  +7.01%  '

# Row 34
$ws.Range("D34").Value = '12.08'
$ws.Range("E34").Value = '  -1.07%  '

# Row 35
$ws.Range("E35").Value = '  -1.66%  '

# Row 36
$ws.Range("D36").Value = '63.75'
$ws.Range("E36").Value = '  -3.70%  '

# Row 37
$ws.Range("D37").Value = '42.24'
$ws.Range("E37").Value = '  -5.54%  '

# Row 38
$ws.Range("D38").Value = '0.406'
$ws.Range("E38").Value = '  +3.94%  '

# Row 39
$ws.Range("E39").Value = '  +0.10%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0765'
$ws.Range("E40").Value = '  -0.22%  '

# Row 41
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").Value = '3.11'
$ws.Range("E41").Value = '  +8.18%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.147.68'
$ws.Range("E42").Value = '  +13.29%  '

# Row 43
$ws.Range("D43").Value = '0.133'
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  +9.77%  '

# Row 45
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.14%  '

# Row 46
$ws.Range("D46").Value = '0.0413'
$ws.Range("E46").Value = '  -1.91%  '

# Row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = '143.49'
$ws.Range("E47").Value = '  +0.53%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.130'
$ws.Range("E48").Value = '  -3.35%  '

# Row 49
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.09'
$ws.Range("E49").Value = '  -1.50%  '

# Row 50
$ws.Range("D50").Value = '8.53'
$ws.Range("E50").Value = '  -2.76%  '

# Row 51
$ws.Range("D51").Value = '2.54'
$ws.Range("E51").Value = '  -2.28%  '
